# CryCompanywiseStockReport_1.xlsx - stock report correction.
#
# A handful of line items had their sold Qty (col F) and Stock Value
# (col G = Unit Cost col D * Qty col F) corrected downward (mostly by a
# few units, some zeroed out). A few pairs of adjacent rows for the same
# product had their Code/MRP/Qty/Value (cols B/E/F/G) swapped back to the
# other order. Every per-company "Sub Total:" (col B) below an affected
# block, plus the two grand-total rows at the bottom, were recalculated
# to match (Sub Total = SUM of col G for that company's rows; the overall
# Sub Total/Grand Total = SUM of all the per-company Sub Totals).
#
# All of the resulting values are applied directly below, cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F71").Value = 340
$ws.Range("G71").Value = 21658

$ws.Range("F77").Value = 266
$ws.Range("G77").Value = 12432.84

$ws.Range("F86").Value = 73
$ws.Range("G86").Value = 9159.309999999999

$ws.Range("B90").Value = 190299.66

$ws.Range("F102").Value = 6
$ws.Range("G102").Value = 296.88

$ws.Range("B104").Value = 264.86

$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 97.2

$ws.Range("F115").Value = 217
$ws.Range("G115").Value = 21007.77

$ws.Range("B117").Value = 14962.85

$ws.Range("F141").Value = 48
$ws.Range("G141").Value = 2569.44

$ws.Range("B142").Value = 3081.12

$ws.Range("F144").Value = 1130
$ws.Range("G144").Value = 9548.5

$ws.Range("B147").Value = 16748.94

$ws.Range("F151").Value = 95
$ws.Range("G151").Value = 8253.6

$ws.Range("B156").Value = 33551.69

$ws.Range("F167").Value = 18
$ws.Range("G167").Value = 5166.18

$ws.Range("B175").Value = 31524.09

$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2

$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 66.59999999999999

$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6

$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18

$ws.Range("F225").Value = 76
$ws.Range("G225").Value = 8681.48

$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 66
$ws.Range("G227").Value = 9522.48

$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32

$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12

$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 119
$ws.Range("G233").Value = 5669.16

$ws.Range("F251").Value = 3
$ws.Range("G251").Value = 734.25

$ws.Range("F255").Value = 588
$ws.Range("G255").Value = 100742.04

$ws.Range("B260").Value = 198694.45

$ws.Range("F277").Value = 1
$ws.Range("G277").Value = 21.25

$ws.Range("F280").Value = 139
$ws.Range("G280").Value = 23510.46

$ws.Range("F282").Value = 4
$ws.Range("G282").Value = 214.8

$ws.Range("F291").Value = 116
$ws.Range("G291").Value = 4989.16

$ws.Range("F293").Value = 45
$ws.Range("G293").Value = 3164.4

$ws.Range("F294").Value = 35
$ws.Range("G294").Value = 2497.6

$ws.Range("F296").Value = 68
$ws.Range("G296").Value = 1441.6

$ws.Range("F302").Value = 64
$ws.Range("G302").Value = 13496.96

$ws.Range("F303").Value = 37
$ws.Range("G303").Value = 7802.93

$ws.Range("B304").Value = 184436.55

$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 0

$ws.Range("B315").Value = 2348.48

$ws.Range("F323").Value = 39
$ws.Range("G323").Value = 4116.06

$ws.Range("B330").Value = 29577.58

$ws.Range("F334").Value = 196
$ws.Range("G334").Value = 10156.72

$ws.Range("F336").Value = 20
$ws.Range("G336").Value = 873

$ws.Range("B346").Value = 27173.57

$ws.Range("F353").Value = 14
$ws.Range("G353").Value = 1920.66

$ws.Range("F355").Value = 13
$ws.Range("G355").Value = 2124.85

$ws.Range("F356").Value = 13
$ws.Range("G356").Value = 1820.13

$ws.Range("B358").Value = 35974.82

$ws.Range("B366").Value = 53263
$ws.Range("E366").Value = 15.29
$ws.Range("F366").Value = -309
$ws.Range("G366").Value = -3958.29

$ws.Range("B367").Value = 65066
$ws.Range("E367").Value = 13.61
$ws.Range("F367").Value = 90
$ws.Range("G367").Value = 1152.9

$ws.Range("B375").Value = 64927
$ws.Range("E375").Value = 17.26
$ws.Range("F375").Value = 106
$ws.Range("G375").Value = 1719.32

$ws.Range("B376").Value = 45718
$ws.Range("E376").Value = 19.38
$ws.Range("F376").Value = -294
$ws.Range("G376").Value = -4768.68

$ws.Range("B380").Value = 45709
$ws.Range("E380").Value = 15.69
$ws.Range("F380").Value = -300
$ws.Range("G380").Value = -3945

$ws.Range("B381").Value = 64925
$ws.Range("E381").Value = 13.97
$ws.Range("F381").Value = 111
$ws.Range("G381").Value = 1459.65

$ws.Range("F399").Value = 0
$ws.Range("G399").Value = 0

$ws.Range("F406").Value = 0
$ws.Range("G406").Value = 0

$ws.Range("B411").Value = 7768.5

$ws.Range("F430").Value = 0
$ws.Range("G430").Value = 0

$ws.Range("F434").Value = 21
$ws.Range("G434").Value = 685.4400000000001

$ws.Range("B435").Value = 830.92

$ws.Range("F450").Value = 12
$ws.Range("G450").Value = 1664.88

$ws.Range("F454").Value = 51
$ws.Range("G454").Value = 1741.65

$ws.Range("B460").Value = 14070.17

$ws.Range("F491").Value = 23
$ws.Range("G491").Value = 4094.46

$ws.Range("B493").Value = 12319.11

$ws.Range("F509").Value = 234
$ws.Range("G509").Value = 18808.92

$ws.Range("B510").Value = 24941.38

$ws.Range("F542").Value = 52
$ws.Range("G542").Value = 6735.56

$ws.Range("B547").Value = 23546.12

$ws.Range("F549").Value = 27
$ws.Range("G549").Value = 1292.22

$ws.Range("F550").Value = 1
$ws.Range("G550").Value = 81.56

$ws.Range("F552").Value = 18
$ws.Range("G552").Value = 1832.22

$ws.Range("F555").Value = 28
$ws.Range("G555").Value = 1947.68

$ws.Range("B560").Value = 5727.98

$ws.Range("F564").Value = 1
$ws.Range("G564").Value = 163.89

$ws.Range("F568").Value = 1
$ws.Range("G568").Value = 157.17

$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4

$ws.Range("F577").Value = 72
$ws.Range("G577").Value = 3095.28

$ws.Range("F580").Value = 62
$ws.Range("G580").Value = 3533.38

$ws.Range("F581").Value = 5
$ws.Range("G581").Value = 1209

$ws.Range("B583").Value = 22038.78

$ws.Range("F599").Value = 1837
$ws.Range("G599").Value = 299633.07

$ws.Range("F601").Value = 434
$ws.Range("G601").Value = 122765.58

$ws.Range("F602").Value = 338
$ws.Range("G602").Value = 48891.7

$ws.Range("B606").Value = 472138.4

$ws.Range("F610").Value = 15
$ws.Range("G610").Value = 614.85

$ws.Range("B618").Value = 44738.42

$ws.Range("B619").Value = 1859017.41

$ws.Range("B620").Value = 1859017.41
